# Update the two-digit multiplication practice sheet: replace each
# equation with the newly generated one. Every "old" string is unique
# in the document, so a literal (non-wildcard) Find/Replace covering
# the whole Content range (wrap mode 1 = wdFindContinue) is enough.
#
# Note on ordering: "17×80=" is both an old value (row 2, last cell)
# and a new value (row 3, last cell, replacing "15×90="). The pairs
# below are ordered so the "17×80=" -> "80×71=" replacement runs
# *before* the "15×90=" -> "17×80=" replacement creates a fresh
# "17×80=" in the document - otherwise that freshly-written text would
# get caught by the earlier rule on a later pass.

$d = $word.ActiveDocument

$pairs = @(
    @("90×71=", "11×26="),
    @("55×65=", "87×97="),
    @("83×49=", "32×93="),
    @("94×49=", "44×95="),
    @("32×79=", "82×65="),
    @("33×45=", "35×25="),
    @("34×31=", "50×42="),
    @("75×79=", "92×96="),
    @("15×16=", "42×95="),
    @("17×80=", "80×71="),
    @("23×85=", "98×39="),
    @("44×66=", "16×97="),
    @("38×13=", "50×87="),
    @("70×99=", "91×82="),
    @("15×90=", "17×80="),
    @("28×48=", "22×34="),
    @("41×44=", "43×51="),
    @("34×14=", "41×92="),
    @("70×19=", "50×53="),
    @("85×55=", "57×73="),
    @("63×13=", "57×11="),
    @("57×25=", "34×85="),
    @("88×64=", "40×40="),
    @("41×82=", "42×44="),
    @("37×55=", "94×53=")
)

$count = 0
foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if ($found) {
        $count = $count + 1
    } else {
        Write-Output "MISSING: $old"
    }
}

Write-Output "Replaced $count of $($pairs.Count) equations"
